# Update TPM-derived NATMI edge statistics for Fgf2-Sdc2 (YoungD7) per new TPM run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.313179
$ws.Range("H2").Value = 0.939537
$ws.Range("I2").Value = 0.02707464596575709
$ws.Range("J2").Value = 0.0270746459657571
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.971099
$ws.Range("N2").Value = 2.913297
$ws.Range("O2").Value = 0.007882574716876797
$ws.Range("P2").Value = 0.007882574716876797
$ws.Range("Q2").Value = 0.304127813721
$ws.Range("R2").Value = 2.737150323489
$ws.Range("S2").Value = 0.0002134179197580672
$ws.Range("T2").Value = 0.0002134179197580673

# Row 3
$ws.Range("G3").Value = 0.313179
$ws.Range("H3").Value = 0.939537
$ws.Range("I3").Value = 0.02707464596575709
$ws.Range("J3").Value = 0.0270746459657571
$ws.Range("O3").Value = 0.6966643430097871
$ws.Range("P3").Value = 0.696664343009787
$ws.Range("Q3").Value = 26.878908370296
$ws.Range("R3").Value = 241.910175332664
$ws.Range("S3").Value = 0.01886194044395675
$ws.Range("T3").Value = 0.01886194044395675

# Row 4
$ws.Range("G4").Value = 0.313179
$ws.Range("H4").Value = 0.939537
$ws.Range("I4").Value = 0.02707464596575709
$ws.Range("J4").Value = 0.0270746459657571
$ws.Range("M4").Value = 36.24916566666667
$ws.Range("N4").Value = 108.747497
$ws.Range("O4").Value = 0.294240604502677
$ws.Range("P4").Value = 0.294240604502677
$ws.Range("Q4").Value = 11.352477454321
$ws.Range("R4").Value = 102.172297088889
$ws.Range("S4").Value = 0.007966460195660333
$ws.Range("T4").Value = 0.007966460195660333

# Row 5
$ws.Range("G5").Value = 0.313179
$ws.Range("H5").Value = 0.939537
$ws.Range("I5").Value = 0.02707464596575709
$ws.Range("J5").Value = 0.0270746459657571
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.149372
$ws.Range("N5").Value = 0.448116
$ws.Range("O5").Value = 0.001212477770659141
$ws.Range("P5").Value = 0.001212477770659141
$ws.Range("Q5").Value = 0.04678017358799999
$ws.Range("R5").Value = 0.4210215622919999
$ws.Range("S5").Value = 0.00003282740638194666
$ws.Range("T5").Value = 0.00003282740638194666

# Row 6
$ws.Range("I6").Value = 0.6982806158817221
$ws.Range("J6").Value = 0.6982806158817222
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.971099
$ws.Range("N6").Value = 2.913297
$ws.Range("O6").Value = 0.007882574716876797
$ws.Range("P6").Value = 0.007882574716876797
$ws.Range("Q6").Value = 7.843742715618667
$ws.Range("R6").Value = 70.593684440568
$ws.Range("S6").Value = 0.005504249128034421
$ws.Range("T6").Value = 0.005504249128034422

# Row 7
$ws.Range("I7").Value = 0.6982806158817221
$ws.Range("J7").Value = 0.6982806158817222
$ws.Range("O7").Value = 0.6966643430097871
$ws.Range("P7").Value = 0.696664343009787
$ws.Range("S7").Value = 0.4864672064997094
$ws.Range("T7").Value = 0.4864672064997094

# Row 8
$ws.Range("I8").Value = 0.6982806158817221
$ws.Range("J8").Value = 0.6982806158817222
$ws.Range("M8").Value = 36.24916566666667
$ws.Range("N8").Value = 108.747497
$ws.Range("O8").Value = 0.294240604502677
$ws.Range("P8").Value = 0.294240604502677
$ws.Range("Q8").Value = 292.7910842717076
$ws.Range("R8").Value = 2635.119758445368
$ws.Range("S8").Value = 0.2054625105295395
$ws.Range("T8").Value = 0.2054625105295395

# Row 9
$ws.Range("I9").Value = 0.6982806158817221
$ws.Range("J9").Value = 0.6982806158817222
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.149372
$ws.Range("N9").Value = 0.448116
$ws.Range("O9").Value = 0.001212477770659141
$ws.Range("P9").Value = 0.001212477770659141
$ws.Range("Q9").Value = 1.206504730122667
$ws.Range("R9").Value = 10.858542571104
$ws.Range("S9").Value = 0.0008466497244387621
$ws.Range("T9").Value = 0.0008466497244387622

# Row 10
$ws.Range("G10").Value = 2.897745666666667
$ws.Range("H10").Value = 8.693237
$ws.Range("I10").Value = 0.2505130868410934
$ws.Range("J10").Value = 0.2505130868410934
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.971099
$ws.Range("N10").Value = 2.913297
$ws.Range("O10").Value = 0.007882574716876797
$ws.Range("P10").Value = 0.007882574716876797
$ws.Range("Q10").Value = 2.813997919154334
$ws.Range("R10").Value = 25.325981272389
$ws.Range("S10").Value = 0.001974688124580364
$ws.Range("T10").Value = 0.001974688124580364

# Row 11
$ws.Range("G11").Value = 2.897745666666667
$ws.Range("H11").Value = 8.693237
$ws.Range("I11").Value = 0.2505130868410934
$ws.Range("J11").Value = 0.2505130868410934
$ws.Range("O11").Value = 0.6966643430097871
$ws.Range("P11").Value = 0.696664343009787
$ws.Range("Q11").Value = 248.7019891332293
$ws.Range("R11").Value = 2238.317902199064
$ws.Range("S11").Value = 0.174523535059504
$ws.Range("T11").Value = 0.174523535059504

# Row 12
$ws.Range("G12").Value = 2.897745666666667
$ws.Range("H12").Value = 8.693237
$ws.Range("I12").Value = 0.2505130868410934
$ws.Range("J12").Value = 0.2505130868410934
$ws.Range("M12").Value = 36.24916566666667
$ws.Range("N12").Value = 108.747497
$ws.Range("O12").Value = 0.294240604502677
$ws.Range("P12").Value = 0.294240604502677
$ws.Range("Q12").Value = 105.0408627308655
$ws.Range("R12").Value = 945.3677645777891
$ws.Range("S12").Value = 0.07371112210795494
$ws.Range("T12").Value = 0.07371112210795493

# Row 13
$ws.Range("G13").Value = 2.897745666666667
$ws.Range("H13").Value = 8.693237
$ws.Range("I13").Value = 0.2505130868410934
$ws.Range("J13").Value = 0.2505130868410934
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.149372
$ws.Range("N13").Value = 0.448116
$ws.Range("O13").Value = 0.001212477770659141
$ws.Range("P13").Value = 0.001212477770659141
$ws.Range("Q13").Value = 0.4328420657213333
$ws.Range("R13").Value = 3.895578591492
$ws.Range("S13").Value = 0.0003037415490540286
$ws.Range("T13").Value = 0.0003037415490540286

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.2791366666666666
$ws.Range("H14").Value = 0.83741
$ws.Range("I14").Value = 0.02413165131142748
$ws.Range("J14").Value = 0.02413165131142749
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.971099
$ws.Range("N14").Value = 2.913297
$ws.Range("O14").Value = 0.007882574716876797
$ws.Range("P14").Value = 0.007882574716876797
$ws.Range("Q14").Value = 0.2710693378633333
$ws.Range("R14").Value = 2.43962404077
$ws.Range("S14").Value = 0.0001902195445039451
$ws.Range("T14").Value = 0.0001902195445039451

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.2791366666666666
$ws.Range("H15").Value = 0.83741
$ws.Range("I15").Value = 0.02413165131142748
$ws.Range("J15").Value = 0.02413165131142749
$ws.Range("O15").Value = 0.6966643430097871
$ws.Range("P15").Value = 0.696664343009787
$ws.Range("Q15").Value = 23.95719025261333
$ws.Range("R15").Value = 215.61471227352
$ws.Range("S15").Value = 0.0168116610066169
$ws.Range("T15").Value = 0.0168116610066169

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.2791366666666666
$ws.Range("H16").Value = 0.83741
$ws.Range("I16").Value = 0.02413165131142748
$ws.Range("J16").Value = 0.02413165131142749
$ws.Range("M16").Value = 36.24916566666667
$ws.Range("N16").Value = 108.747497
$ws.Range("O16").Value = 0.294240604502677
$ws.Range("P16").Value = 0.294240604502677
$ws.Range("Q16").Value = 10.11847127364111
$ws.Range("R16").Value = 91.06624146277001
$ws.Range("S16").Value = 0.007100511669522242
$ws.Range("T16").Value = 0.007100511669522243

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.2791366666666666
$ws.Range("H17").Value = 0.83741
$ws.Range("I17").Value = 0.02413165131142748
$ws.Range("J17").Value = 0.02413165131142749
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.149372
$ws.Range("N17").Value = 0.448116
$ws.Range("O17").Value = 0.001212477770659141
$ws.Range("P17").Value = 0.001212477770659141
$ws.Range("Q17").Value = 0.04169520217333332
$ws.Range("R17").Value = 0.37525681956
$ws.Range("S17").Value = 0.00002925909078440333
$ws.Range("T17").Value = 0.00002925909078440333
